# Weekly fruit/hortaliza update: a new weekly price record is inserted
# as row 34 (Macroferia Regional de Talca - Espárragos, Provincia de
# Linares, week of 2021-11-08), pushing the previously-existing rows
# 34-44 down to rows 35-45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 34, shifting rows 34:44 -> 35:45
$ws.Rows("34:34").Insert()

# Populate the newly inserted row 34 with this week's record
$ws.Range("A34").Value = 5
$ws.Range("B34").Value = "Macroferia Regional de Talca"
$ws.Range("C34").Value = "Maule"
$ws.Range("D34").Value = 44508
$ws.Range("E34").Value = 7
$ws.Range("F34").Value = 300000000
$ws.Range("G34").Value = "Espárragos"
$ws.Range("H34").Value = "Verde"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 800
$ws.Range("L34").Value = 800
$ws.Range("M34").Value = 800
$ws.Range("N34").Value = "$/kilo"
$ws.Range("O34").Value = "Provincia de Linares"
$ws.Range("P34").Value = 800
$ws.Range("Q34").Value = 1
$ws.Range("R34").Value = "Hortaliza"
